# "pówer bi com esboços de visualização"
# Re-label the age-range lookup table (d_faixas_idade) so the two
# open-ended buckets read as "<= 20 anos" / "> 50 " instead of
# "Menos de 20 anos" / "50 anos ou mais", and sketch out a new
# (still-empty) row below the table for further visualization work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("d_faixas_idade")

# Rename the open-ended buckets.
$ws.Range("B2").Value = "<= 20 anos"
$ws.Range("B6").Value = "> 50 "

# Sketch a new row (B13) below the table: give it the same underlined
# style used elsewhere in the workbook for draft/placeholder cells,
# without putting any value in it yet.
$ws.Range("B13").Font.Underline = $true

# Leave the selection where the new sketch row is.
$ws.Range("B13").Select() | Out-Null
